# Update cryptos list values (Price and Volume(1h) columns) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.024.58'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.677.66'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.89'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.32'
$ws.Range('E9').Value = '  +5.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0622'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0888'
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.914.53'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.677.47'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.47'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.034.73'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.12'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '235.53'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0738'
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.26'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('E24').Value = '  -4.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.70'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.44'
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.536.46'
$ws.Range('E33').Value = '  +5.23%  '
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.918'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('E40').Value = '  +6.39%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '67.94'
$ws.Range('E42').Value = '  +3.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.58'
$ws.Range('E43').Value = '  -2.57%  '
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.819.32'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.782'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.34'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.08'
$ws.Range('E50').Value = '  +7.27%  '
$ws.Range('E51').Value = '  -0.32%  '
